$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" list in E16:E52 was sorted descending (2104 -> 1804).
# It is being re-sorted ascending (1804 -> 2104). Capture current labels,
# reverse them, and write them back so the underlying shared-string table
# and the visible labels match the new ordering.
$periods = @()
for ($r = 16; $r -le 52; $r++) {
    $periods += $ws.Cells.Item($r, 5).Value2
}

$count = $periods.Length
for ($i = 0; $i -lt $count; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$count - 1 - $i]
}

# The "Valor Mora" amounts travel with their period label: period 1804's
# amount (35200) now lands on row 16, and period 2104's amount (30506)
# now lands on row 52.
$ws.Cells.Item(16, 6).Value = 35200
$ws.Cells.Item(52, 6).Value = 30506
